{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Source change (from the OOXML diff): the single paragraph's runs\n// (\"TEST \" + \"3\") are replaced with a large \"Ideal Job ...\" passage\n// ending in \"... the learning styles \", followed by a short run \"tes\"\n// that is wrapped in proofing-error markers\n// (<w:proofErr w:type=\"spellStart\"/> ... <w:proofErr w:type=\"spellEnd\"/>),\n// Word's markup for flagging \"tes\" as a misspelling. The explicit\n// w:lang=\"en-US\" formatting that lived on the paragraph mark and on both\n// runs is also dropped.\n//\n// The Office.js object model has no direct property for w:proofErr (it\n// is a proofing artifact, not a content API), so we build the exact\n// OOXML for the new paragraph ourselves and insert it as \"flat OPC\"\n// package XML via Body.insertOoxml(..., Word.InsertLocation.replace),\n// which swaps out the whole body content in one shot.\n\nconst bodyText = \"Ideal Job Find an advertisement showing what you believe to be your ideal job. This may require several years of experience, and hence be something that you must work towards, rather than something that you are ready for now, or will be able to fill as soon as you graduate. There are various ways to search for IT jobs, including websites like {seek.com.au}. You should include the following information. \\u2022 The job advertisement itself. Include a link, and a snapshot of it (in case the link expires before the assignment deadline). \\u2022 A description (in your own words) of the position, and particularly what makes this position appealing to you. \\u2022 A description (in your own words) of the skills, qualifications and experience required for the position. \\u2022 A description (in your own words) of the skills, qualifications and experience you currently have. \\u2022 A plan describing how you will obtain the skills, qualifications and experience required for the position, building on those you have now. This need not be greatly detailed, (and will probably change significantly over time anyway), but try to be as specific as you can. Personal Profile There are a number of online tests that are commonly used by employers to get specific information about potential employees. One of the best-known of these is the Myers-Briggs Type Indicator (MBTI) test, which was developed by the mother-and-daughter team of Isabel Myers and Katharine Briggs. Another popular one (and very relevant to students) is about learning styles. There are various other tests available online as well, and while there is no guarantee that any specific test will be necessarily one that you will encounter in your career, it seems highly likely that you will have to do some kind of test like this as part of a recruitment process. Accordingly, you are required to present the following information: \\u2022 The results of an online Myers-Briggs test. www.16personalities.com \\u2022 The results of an online learning style test. \\u2022 The results of one further online test of your choosing. The third test should be distinct from both the Myers-Briggs test and the learning styles \";\n\nconst paragraphOoxml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' +\n    '<w:p w14:paraId=\"33AFCCDF\" w14:textId=\"62A5D320\" w:rsidR=\"0097470C\" w:rsidRPr=\"00C24C45\" w:rsidRDefault=\"00C24C45\" w:rsidP=\"006302A1\">' +\n      '<w:r><w:t xml:space=\"preserve\">' + bodyText + '</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r w:rsidR=\"00D15A51\"><w:t>tes</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>' +\n  '</w:body>' +\n  '</w:document>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>\\n' +\n  '<?mso-application progid=\"Word.Document\"?>\\n' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' + paragraphOoxml + '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nbody.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Source change (from the OOXML diff): the document's single paragraph had\n# two runs, \"TEST \" and \"3\". Both are replaced by a large \"Ideal Job ...\"\n# passage (ending \"... the learning styles \") followed by a short run\n# \"tes\" that Word has flagged as a misspelling, i.e. wrapped in\n# <w:proofErr w:type=\"spellStart\"/> ... <w:proofErr w:type=\"spellEnd\"/>.\n# The explicit w:lang=\"en-US\" formatting that was on the paragraph mark\n# and on both runs is also removed.\n#\n# The Word object model has no direct property for w:proofErr (it's a\n# proofing-pass artifact, not a content API), so the new paragraph is\n# built as literal OOXML and dropped in with Range.InsertXML, which is\n# COM's equivalent of Office.js's Range.insertOoxml and accepts the same\n# \"flat OPC\" package XML.\n\n$d = $word.ActiveDocument\n\n$bodyText = \"Ideal Job Find an advertisement showing what you believe to be your ideal job. This may require several years of experience, and hence be something that you must work towards, rather than something that you are ready for now, or will be able to fill as soon as you graduate. There are various ways to search for IT jobs, including websites like {seek.com.au}. You should include the following information. \u2022 The job advertisement itself. Include a link, and a snapshot of it (in case the link expires before the assignment deadline). \u2022 A description (in your own words) of the position, and particularly what makes this position appealing to you. \u2022 A description (in your own words) of the skills, qualifications and experience required for the position. \u2022 A description (in your own words) of the skills, qualifications and experience you currently have. \u2022 A plan describing how you will obtain the skills, qualifications and experience required for the position, building on those you have now. This need not be greatly detailed, (and will probably change significantly over time anyway), but try to be as specific as you can. Personal Profile There are a number of online tests that are commonly used by employers to get specific information about potential employees. One of the best-known of these is the Myers-Briggs Type Indicator (MBTI) test, which was developed by the mother-and-daughter team of Isabel Myers and Katharine Briggs. Another popular one (and very relevant to students) is about learning styles. There are various other tests available online as well, and while there is no guarantee that any specific test will be necessarily one that you will encounter in your career, it seems highly likely that you will have to do some kind of test like this as part of a recruitment process. Accordingly, you are required to present the following information: \u2022 The results of an online Myers-Briggs test. www.16personalities.com \u2022 The results of an online learning style test. \u2022 The results of one further online test of your choosing. The third test should be distinct from both the Myers-Briggs test and the learning styles \"\n\n$paragraphOoxml = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"33AFCCDF\" w14:textId=\"62A5D320\" w:rsidR=\"0097470C\" w:rsidRPr=\"00C24C45\" w:rsidRDefault=\"00C24C45\" w:rsidP=\"006302A1\"><w:r><w:t xml:space=\"preserve\">' + $bodyText + '</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r w:rsidR=\"00D15A51\"><w:t>tes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p></w:body></w:document>'\n\n$flatOpc = '<?xml version=\"1.0\" standalone=\"yes\"?>' + \"`n\" + '<?mso-application progid=\"Word.Document\"?>' + \"`n\" + '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + $paragraphOoxml + '</pkg:xmlData></pkg:part></pkg:package>'\n\n# Replace the whole story (all paragraphs + the final paragraph mark) with\n# the freshly built OOXML in one operation.\n$r = $d.Content\n$r.InsertXML($flatOpc)\n"}
